# Update "want-to-go" counts (column F) across all sheets to match
# the regenerated gh-pages data output at commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 585
$ws.Cells.Item(3, 6).Value = 10510
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(6, 6).Value = 6995
$ws.Cells.Item(7, 6).Value = 658
$ws.Cells.Item(8, 6).Value = 141
$ws.Cells.Item(9, 6).Value = 12534
$ws.Cells.Item(10, 6).Value = 12534
$ws.Cells.Item(11, 6).Value = 12896
$ws.Cells.Item(12, 6).Value = 1316
$ws.Cells.Item(13, 6).Value = 1290
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(17, 6).Value = 363
$ws.Cells.Item(19, 6).Value = 1431
$ws.Cells.Item(20, 6).Value = 355
$ws.Cells.Item(21, 6).Value = 2016
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(23, 6).Value = 1551
$ws.Cells.Item(24, 6).Value = 881
$ws.Cells.Item(25, 6).Value = 18
$ws.Cells.Item(27, 6).Value = 733
$ws.Cells.Item(28, 6).Value = 3002
$ws.Cells.Item(29, 6).Value = 254
$ws.Cells.Item(30, 6).Value = 2049
$ws.Cells.Item(31, 6).Value = 109
$ws.Cells.Item(33, 6).Value = 1000
$ws.Cells.Item(34, 6).Value = 143
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(36, 6).Value = 98
$ws.Cells.Item(37, 6).Value = 3746
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(40, 6).Value = 130
$ws.Cells.Item(43, 6).Value = 1470
$ws.Cells.Item(44, 6).Value = 35
$ws.Cells.Item(47, 6).Value = 37
$ws.Cells.Item(48, 6).Value = 31
$ws.Cells.Item(49, 6).Value = 4295
$ws.Cells.Item(50, 6).Value = 188

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(4, 6).Value = 36
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(7, 6).Value = 38
$ws.Cells.Item(8, 6).Value = 80
$ws.Cells.Item(25, 6).Value = 108
$ws.Cells.Item(26, 6).Value = 49

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 6516
$ws.Cells.Item(3, 6).Value = 30

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 585
$ws.Cells.Item(3, 6).Value = 10510
$ws.Cells.Item(6, 6).Value = 6995
$ws.Cells.Item(7, 6).Value = 658
$ws.Cells.Item(10, 6).Value = 12534
$ws.Cells.Item(11, 6).Value = 12896
$ws.Cells.Item(12, 6).Value = 38
$ws.Cells.Item(13, 6).Value = 1316
$ws.Cells.Item(14, 6).Value = 1290
$ws.Cells.Item(15, 6).Value = 911
$ws.Cells.Item(16, 6).Value = 363
$ws.Cells.Item(17, 6).Value = 80
$ws.Cells.Item(18, 6).Value = 194
$ws.Cells.Item(21, 6).Value = 2016
$ws.Cells.Item(23, 6).Value = 1551
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(27, 6).Value = 3002
$ws.Cells.Item(28, 6).Value = 254
$ws.Cells.Item(29, 6).Value = 2049
$ws.Cells.Item(30, 6).Value = 109
$ws.Cells.Item(32, 6).Value = 1683
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(36, 6).Value = 143
$ws.Cells.Item(38, 6).Value = 98
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(41, 6).Value = 49
$ws.Cells.Item(42, 6).Value = 275
$ws.Cells.Item(44, 6).Value = 646
$ws.Cells.Item(45, 6).Value = 1470
$ws.Cells.Item(48, 6).Value = 31
$ws.Cells.Item(49, 6).Value = 4295
$ws.Cells.Item(50, 6).Value = 188

